$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell formatting from G6 down to G7:G8 so the new
# date cells reuse the existing style (numFmtId 22) instead of minting
# a brand-new style entry.
$ws.Cells.Item(6, 7).Copy()
$ws.Range($ws.Cells.Item(7, 7), $ws.Cells.Item(8, 7)).PasteSpecial(-4122)

# Row 7
$ws.Cells.Item(7, 1).Value = 9969.64
$ws.Cells.Item(7, 2).Value = 9950.73
$ws.Cells.Item(7, 3).Value = 78.05
$ws.Cells.Item(7, 4).Value = 78.2
$ws.Cells.Item(7, 5).Value = $false
$ws.Cells.Item(7, 6).Value = 0.19
$ws.Cells.Item(7, 7).Value = 42613.766331018516
$ws.Cells.Item(7, 8).Value = $true

# Row 8
$ws.Cells.Item(8, 1).Value = 10051.39
$ws.Cells.Item(8, 2).Value = 9969.64
$ws.Cells.Item(8, 3).Value = 77.74
$ws.Cells.Item(8, 4).Value = 78.38
$ws.Cells.Item(8, 5).Value = $false
$ws.Cells.Item(8, 6).Value = 0.82
$ws.Cells.Item(8, 7).Value = 42614.674386574072
$ws.Cells.Item(8, 8).Value = $true
